# Append 5 new daily COVID summary rows (2021-10-08 .. 2021-10-12) to the
# "covid_totals" sheet, continuing directly after the existing last row (422).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 423; Date = "2021-10-08"; AreaType = "overview"; AreaCode = "K02000001"; AreaName = "United Kingdom"; CumCases = 8081300; NewCases = 36060; NewDeaths = 127;  CumDeaths = 137541 },
    @{ Row = 424; Date = "2021-10-09"; AreaType = "overview"; AreaCode = "K02000001"; AreaName = "United Kingdom"; CumCases = 8120713; NewCases = 34950; NewDeaths = 133;  CumDeaths = 137697 },
    @{ Row = 425; Date = "2021-10-10"; AreaType = "overview"; AreaCode = "K02000001"; AreaName = "United Kingdom"; CumCases = 8154306; NewCases = 34574; NewDeaths = 38;   CumDeaths = 137735 },
    @{ Row = 426; Date = "2021-10-11"; AreaType = "overview"; AreaCode = "K02000001"; AreaName = "United Kingdom"; CumCases = 8193769; NewCases = 40224; NewDeaths = 28;   CumDeaths = 137763 },
    @{ Row = 427; Date = "2021-10-12"; AreaType = "overview"; AreaCode = "K02000001"; AreaName = "United Kingdom"; CumCases = 8231437; NewCases = 38520; NewDeaths = 181;  CumDeaths = 137944 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Column A holds a date formatted as "yyyy-mm-dd" TEXT (not an Excel
    # date serial) in every other row of this column, so force text entry
    # via a temporary "@" number format, then restore the default "Normal"
    # style so no stray style index is left behind on the cell.
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Date
    $cellA.Style = "Normal"

    $ws.Cells.Item($rowNum, 2).Value = $r.AreaType
    $ws.Cells.Item($rowNum, 3).Value = $r.AreaCode
    $ws.Cells.Item($rowNum, 4).Value = $r.AreaName
    $ws.Cells.Item($rowNum, 5).Value = $r.CumCases
    $ws.Cells.Item($rowNum, 6).Value = $r.NewCases
    $ws.Cells.Item($rowNum, 7).Value = $r.NewDeaths
    $ws.Cells.Item($rowNum, 8).Value = $r.CumDeaths
}
